$wb = $excel.ActiveWorkbook

$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# --- "VENTAS POR GRUPO" sheet: per-client monthly breakdown values ---
$wsVentasPorGrupo.Range("M4").Value = 403.45
$wsVentasPorGrupo.Range("D5").Value = 475.2
$wsVentasPorGrupo.Range("L5").Value = 556.8099999999999
$wsVentasPorGrupo.Range("M5").Value = 151.34
$wsVentasPorGrupo.Range("L12").Value = 886.88
$wsVentasPorGrupo.Range("L19").Value = 1140.32
$wsVentasPorGrupo.Range("E24").Value = 222.87
$wsVentasPorGrupo.Range("M24").Value = 1472.64
$wsVentasPorGrupo.Range("I36").Value = 542.7
$wsVentasPorGrupo.Range("P36").Value = 406.65
$wsVentasPorGrupo.Range("I37").Value = 2168.1
$wsVentasPorGrupo.Range("D50").Value = 91.58

# Row 56 "N de 54" completion counters recomputed for the columns that
# gained a new non-zero value above.
$wsVentasPorGrupo.Range("D56").Value = "2 de 54"
$wsVentasPorGrupo.Range("E56").Value = "1 de 54"
$wsVentasPorGrupo.Range("I56").Value = "3 de 54"
$wsVentasPorGrupo.Range("L56").Value = "3 de 54"
$wsVentasPorGrupo.Range("M56").Value = "6 de 54"
$wsVentasPorGrupo.Range("P56").Value = "1 de 54"

# --- "VENTA MENSUAL" sheet: monthly totals per client ---
$wsVentaMensual.Range("F4").Value = 403.45
$wsVentaMensual.Range("F5").Value = 1183.35
$wsVentaMensual.Range("F12").Value = 886.88
$wsVentaMensual.Range("F19").Value = 1140.32
$wsVentaMensual.Range("F24").Value = 1695.51
$wsVentaMensual.Range("F36").Value = 949.35
$wsVentaMensual.Range("F37").Value = 2168.1
$wsVentaMensual.Range("F50").Value = 91.58
$wsVentaMensual.Range("F60").Value = 12831.54
